$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("API")
$ws2 = $wb.Worksheets.Item("CAU")

# Rename sheets: API -> POP1, CAU -> POP2
$ws1.Name = "POP1"
$ws2.Name = "POP2"

# Update recalculated statistics on sheet1 (formerly API)
$ws1.Range("C50").Value = 1.29846553102691
$ws1.Range("D50").Value = 2.83890139354767
$ws1.Range("E50").Value = 0.0000417788198375411
$ws1.Range("C52").Value = 1.15290397562118
$ws1.Range("D52").Value = 2.59109920537086
$ws1.Range("E52").Value = 0.00149642732357132
$ws1.Range("C59").Value = 1.17691694241317
$ws1.Range("D59").Value = 2.44599959153024
$ws1.Range("E59").Value = 0.00130574243666781
$ws1.Range("C60").Value = 1.15297297041056
$ws1.Range("D60").Value = 2.53894917644285
$ws1.Range("E60").Value = 0.00218675537677159
$ws1.Range("C76").Value = 0.378649954085721
$ws1.Range("D76").Value = 0.872938956201037
$ws1.Range("E76").Value = 0.00617243686350434
$ws1.Range("C80").Value = 0.329374975960078
$ws1.Range("D80").Value = 0.794629029293339
$ws1.Range("E80").Value = 0.000326060509719352
$ws1.Range("C82").Value = 0.421673769947399
$ws1.Range("D82").Value = 0.868080966894022
$ws1.Range("E82").Value = 0.00356564141394195
$ws1.Range("C90").Value = 1.15773701592881
$ws1.Range("D90").Value = 2.42155038231304
$ws1.Range("E90").Value = 0.00183349493082101
$ws1.Range("C98").Value = 0.183431787305971
$ws1.Range("D98").Value = 0.751639499106748
$ws1.Range("E98").Value = 0.00218675537677159
$ws1.Range("C100").Value = 0.382123918381923
$ws1.Range("D100").Value = 0.88624008813449
$ws1.Range("E100").Value = 0.00617243686350434

# Update recalculated statistics on sheet2 (formerly CAU)
$ws2.Range("C20").Value = 1.43471174792119
$ws2.Range("D20").Value = 2.6043780316822
$ws2.Range("E20").Value = 0.0000000081709421276241
$ws2.Range("C55").Value = 0.487205212235643
$ws2.Range("D55").Value = 0.90988358787597
$ws2.Range("E55").Value = 0.00424370571038341
$ws2.Range("C79").Value = 0.603916174646093
$ws2.Range("D79").Value = 0.975496241612133
$ws2.Range("E79").Value = 0.024747308543632
